# Refresh the crypto price/volume snapshot (GitHub Actions scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.296.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.710.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "660.69"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.425"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.707.73"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000320"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +19.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.40"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.401.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.125.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.708.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("E22").Value = "  -3.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "521.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000218"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.47%  "
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "102.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  +15.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "
$ws.Range("E34").Value = "  -0.78%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "656.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.595"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.91%  "
$ws.Range("E43").Value = "  +3.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.493"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.68%  "
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0461"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.26%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.09%  "
